# Playing around with strength analysis:
# Insert a new "0.1% Design Strength" column before the elastic modulus
# columns, populate it with the new per-series values, and remove the
# extra "X" series row (row 13) that is no longer part of the analysis.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column I (old I/J/K -> J/K/L).
$ws.Columns("I:I").Insert()

# New header for the inserted column.
$ws.Range("I1").Value = "0.1% Design Strength"

# New 0.1% design-strength values for the series that have full stats.
$ws.Range("I2").Value = 23.20243803640314
$ws.Range("I3").Value = 33.07472848284672
$ws.Range("I4").Value = 55.52712910140356
$ws.Range("I5").Value = 36.58608608616619
$ws.Range("I6").Value = 13.69604116136211
$ws.Range("I12").Value = 37.36461920746228

# Remove the old "X" series row (row 13) entirely.
$ws.Rows("13:13").Delete()
